$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 151; this shifts existing rows 151-226 down to 152-227
$ws.Rows(151).Insert()

# Populate the newly inserted row 151 with its data
$ws.Range("A151").Value = 10
$ws.Range("B151").Value = "Vega Modelo de Temuco"
$ws.Range("C151").Value = "La Araucanía"
$ws.Range("D151").Value = 44518
$ws.Range("E151").Value = 9
$ws.Range("F151").Value = 100112044
$ws.Range("G151").Value = "Perejil"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 60
$ws.Range("K151").Value = 4000
$ws.Range("L151").Value = 5000
$ws.Range("M151").Value = 4500
$ws.Range("N151").Value = "$/docena de atados (3 kilos)"
$ws.Range("O151").Value = "Provincia de Cautín"
$ws.Range("P151").Value = 1500
$ws.Range("Q151").Value = 3
$ws.Range("R151").Value = "Hortaliza"
